# API over 50% complete
# Update the "hobi api docs" tracking sheet:
#  - mark several previously "NOT STARTED" webmethods as "COMPLETED"
#  - rename two moderator-related rows
#  - add a new "Needs testing" note
#  - add four new tracked API rows (get group moderators / set group admin /
#    add+remove group moderator)
#  - recompute the completion percentages

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Webmethods column (C) for rows 2-17 goes from NOT STARTED -> COMPLETED
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 17; $r++) {
    $ws.Range("C$r").Value = "COMPLETED"
}
$ws.Range("B2:B17").Copy() | Out-Null
$ws.Range("C2:C17").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 2) Rows 18 & 19 (stored procs + webmethods) go from NOT STARTED -> COMPLETED
# ---------------------------------------------------------------------------
$ws.Range("B18").Value = "COMPLETED"
$ws.Range("C18").Value = "COMPLETED"
$ws.Range("B19").Value = "COMPLETED"
$ws.Range("C19").Value = "COMPLETED"

$ws.Range("B2").Copy() | Out-Null
$ws.Range("B18").PasteSpecial(-4122) | Out-Null
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B19").PasteSpecial(-4122) | Out-Null
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C18").PasteSpecial(-4122) | Out-Null
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C19").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 3) Rename the two moderator rows to call out private groups
# ---------------------------------------------------------------------------
$ws.Range("A27").Value = "(mod) accept member - Private group"
$ws.Range("A28").Value = "(mod) decline member - private group"

# ---------------------------------------------------------------------------
# 4) New "Needs testing" callout next to the "edit user image" row
# ---------------------------------------------------------------------------
$ws.Range("D4").Value = "> Needs testing"

# ---------------------------------------------------------------------------
# 5) "edit user password" (row 35) regresses from COMPLETED -> NOT STARTED
# ---------------------------------------------------------------------------
$ws.Range("B35").Value = "NOT STARTED"
$ws.Range("C20").Copy() | Out-Null
$ws.Range("B35").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 6) Insert three extra rows before the old blank separator row (row 36) so
#    there is room for four new tracked API rows, then fill them in.
# ---------------------------------------------------------------------------
$ws.Rows.Item(36).Resize(3).Insert()

# Rows 36-38 are brand new; row 39 is the old blank/thick-bordered separator
# row that got pushed down (it already carries the thick-bottom row format).
$newRows = @(
    @{ Row = 36; Name = "get group moderators"; Status = "NOT STARTED" },
    @{ Row = 37; Name = "set group admin"; Status = "COMPLETED" },
    @{ Row = 38; Name = "(admin) add group moderator"; Status = "COMPLETED" },
    @{ Row = 39; Name = "(admin) Remove group moderator"; Status = "COMPLETED" }
)

foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Range("A$r").Value = $item.Name
    $ws.Range("B$r").Value = $item.Status
    $ws.Range("C$r").Value = $item.Status

    if ($item.Status -eq "COMPLETED") {
        $ws.Range("B2:C2").Copy() | Out-Null
    } else {
        $ws.Range("B20:C20").Copy() | Out-Null
    }
    $ws.Range("B$r`:C$r").PasteSpecial(-4122) | Out-Null
}

# Row 39 needs the thick bottom row height/border restored (the insert above
# keeps it on the row that already had it, but make sure it's explicit).
$ws.Rows.Item(39).RowHeight = 14.7

# ---------------------------------------------------------------------------
# 7) Totals row (now row 40): update formulas for the new 21/38 ratio
# ---------------------------------------------------------------------------
$ws.Range("B40").Formula = "=21/38 * 100"
$ws.Range("C40").Formula = "=(21/38) * 100"
$ws.Range("D40").Formula = "=(B40+C40)/2"

# ---------------------------------------------------------------------------
# 8) Sheet level cosmetics: dimension grows automatically; update column
#    widths, zoom and selection to match.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 31.5
$ws.Columns.Item(4).ColumnWidth = 16.67

$excel.ActiveWindow.Zoom = 70
$ws.Range("E36").Select() | Out-Null
